$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "getspo2ida_parameters"
$ws.Range("B5").Value = "X"
$ws.Range("B6").Value = "X"

$ws.Columns.Item(1).AutoFit() | Out-Null
$ws.Columns.Item(3).AutoFit() | Out-Null
$ws.Columns.Item(5).AutoFit() | Out-Null

$ws.Range("B6").Select()
